$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Hunk 0
$ws.Range("H33").Value = 141.45833
$ws.Range("I33").Value = 91.8421
$ws.Range("J33").Value = 330
$ws.Range("K33").Value = 91.8421
$ws.Range("L33").Value = 330
$ws.Range("M33").Value = 137.1579
$ws.Range("N33").Value = -788
# Hunk 1
$ws.Range("H39").Value = 49.642857
$ws.Range("I39").Value = 28.7
$ws.Range("J39").Value = 102
$ws.Range("K39").Value = 86.09999999999999
$ws.Range("L39").Value = 306
$ws.Range("M39").Value = 209.9
$ws.Range("N39").Value = -898
# Hunk 2
$ws.Range("H95").Value = 15799.625
$ws.Range("J95").Value = 15799.625
$ws.Range("L95").Value = 15799.625
$ws.Range("N95").Value = -21291.625
# Hunk 3
$ws.Range("H141").Value = 3653.6667
$ws.Range("I141").Value = 2235.375
$ws.Range("K141").Value = 6706.125
$ws.Range("M141").Value = -1526.125

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Hunk 4
$ws.Range("H10").Value = 7499.625
$ws.Range("I10").Value = 7799.4
$ws.Range("K10").Value = 7799.4
$ws.Range("M10").Value = -7629.4
# Hunk 5
$ws.Range("H32").Value = 6894.44
$ws.Range("I32").Value = 5075.6113
$ws.Range("J32").Value = 11571.429
$ws.Range("K32").Value = 5075.6113
$ws.Range("L32").Value = 11571.429
$ws.Range("M32").Value = -4788.6113
$ws.Range("N32").Value = -12145.429
# Hunk 6
$ws.Range("H45").Value = 2647
$ws.Range("I45").Value = 1870.5
$ws.Range("K45").Value = 1870.5
$ws.Range("M45").Value = -1493.5
# Hunk 7
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Hunk 8
$ws.Range("H11").Value = 2352.7693
$ws.Range("J11").Value = 2686
$ws.Range("L11").Value = 2686
$ws.Range("N11").Value = -2966
# Hunk 9
$ws.Range("H17").Value = 2300
$ws.Range("J17").Value = 2300
$ws.Range("L17").Value = 2300
$ws.Range("N17").Value = -2644
# Hunk 10
$ws.Range("H86").Value = 3583.3333
# Hunk 11
$ws.Range("H89").Value = 3583.3333
# Hunk 12
$ws.Range("H94").Value = 1441.4348
$ws.Range("I94").Value = 457.9091
$ws.Range("J94").Value = 2343
$ws.Range("K94").Value = 457.9091
$ws.Range("L94").Value = 2343
$ws.Range("M94").Value = -6.909100000000024
$ws.Range("N94").Value = -3245
# Hunk 13
$ws.Range("H107").Value = 997
$ws.Range("I107").Value = 996
$ws.Range("K107").Value = 996
$ws.Range("M107").Value = 924

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Hunk 14
$ws.Range("H16").Value = 100003420
$ws.Range("I16").Value = 100003420
$ws.Range("K16").Value = 100003420
$ws.Range("M16").Value = -100003133
# Hunk 15
$ws.Range("H17").Value = 800
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 800
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -626
$ws.Range("N17").ClearContents()
# Hunk 16
$ws.Range("H22").Value = 100149
$ws.Range("I22").Value = 100163.8
$ws.Range("J22").Value = 100001
$ws.Range("K22").Value = 100163.8
$ws.Range("L22").Value = 100001
$ws.Range("M22").Value = -99813.8
$ws.Range("N22").Value = -100701
# Hunk 17
$ws.Range("H32").Value = 2010
$ws.Range("I32").Value = 2010
$ws.Range("K32").Value = 2010
$ws.Range("M32").Value = -1694
# Hunk 18
$ws.Range("H50").Value = 17023.572
$ws.Range("J50").Value = 23799.8
$ws.Range("L50").Value = 23799.8
$ws.Range("N50").Value = -25049.8
# Hunk 19
$ws.Range("H58").Value = 3753.6924
$ws.Range("I58").Value = 2899.1667
$ws.Range("K58").Value = 2899.1667
$ws.Range("M58").Value = -2696.1667
# Hunk 20
$ws.Range("H113").Value = 100003420
$ws.Range("I113").Value = 100003420
$ws.Range("K113").Value = 100003420
$ws.Range("M113").Value = -100001250
# Hunk 21
$ws.Range("H132").Value = 2584.9312
$ws.Range("I132").Value = 2368.8
$ws.Range("J132").Value = 3065.2222
$ws.Range("K132").Value = 7106.400000000001
$ws.Range("L132").Value = 9195.6666
$ws.Range("M132").Value = -4576.400000000001
$ws.Range("N132").Value = -14255.6666
# Hunk 22
$ws.Range("H136").Value = 3753.6924
$ws.Range("I136").Value = 2899.1667
$ws.Range("K136").Value = 8697.500100000001
$ws.Range("M136").Value = -6147.500100000001

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Hunk 23
$ws.Range("H49").Value = 486
$ws.Range("J49").Value = 430
$ws.Range("L49").Value = 1290
$ws.Range("N49").Value = -1602

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Hunk 24
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31262
# Hunk 25
$ws.Range("H70").Value = 6832.1665
$ws.Range("I70").Value = 5666
$ws.Range("J70").Value = 7998.3335
$ws.Range("K70").Value = 5666
$ws.Range("L70").Value = 7998.3335
$ws.Range("M70").Value = -5396
$ws.Range("N70").Value = -8538.333500000001
# Hunk 26
$ws.Range("H73").Value = 6832.1665
$ws.Range("I73").Value = 5666
$ws.Range("J73").Value = 7998.3335
$ws.Range("K73").Value = 5666
$ws.Range("L73").Value = 7998.3335
$ws.Range("M73").Value = -4730
$ws.Range("N73").Value = -9870.333500000001
# Hunk 27
$ws.Range("H97").Value = 1794.1538
$ws.Range("I97").Value = 979.58826
$ws.Range("K97").Value = 979.58826
$ws.Range("M97").Value = -483.58826
# Hunk 28
$ws.Range("H102").Value = 1074.8125
$ws.Range("I102").Value = 554
$ws.Range("J102").Value = 1942.8334
$ws.Range("K102").Value = 554
$ws.Range("L102").Value = 1942.8334
$ws.Range("M102").Value = 1068
$ws.Range("N102").Value = -5186.8334

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Hunk 29
$ws.Range("H46").Value = 2461.125
$ws.Range("I46").Value = 1465.3334
$ws.Range("K46").Value = 1465.3334
$ws.Range("M46").Value = -1277.3334
# Hunk 30
$ws.Range("H55").Value = 548.8095
$ws.Range("I55").Value = 515.73334
$ws.Range("K55").Value = 515.73334
$ws.Range("M55").Value = -342.73334

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Hunk 31
$ws.Range("H20").Value = 16332.223
$ws.Range("J20").Value = 13374.375
$ws.Range("L20").Value = 13374.375
$ws.Range("N20").Value = -13854.375
# Hunk 32
$ws.Range("H113").Value = 2314
$ws.Range("I113").Value = 776.8
$ws.Range("K113").Value = 2330.4
$ws.Range("M113").Value = -160.3999999999996
# Hunk 33
$ws.Range("H132").Value = 2254.2307
$ws.Range("J132").Value = 2123.5
$ws.Range("L132").Value = 6370.5
$ws.Range("N132").Value = -11430.5
